$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 3 ("汽車" / cars): the old sheet had two badly-shaped rows (a car
# row duplicated into row1/row2, and a stray "boat" row with its own mini
# header baked into the data). The fix folds everything into a single,
# correctly-shaped header + data row pair that matches the other sheets
# (name / capacity / owner / register_date / register_reason /
#  acquire_value / property_category / category / date / legislator_name /
#  legislator_id / source_file / index).
# ---------------------------------------------------------------------------
$carSheet = $wb.Worksheets.Item(3)

# Drop the old third row (the stray boat/header row) entirely so the used
# range shrinks back down to just the header + one data row.
$carSheet.Rows.Item(3).Delete()

# Make sure the new header cells (H1:N1) inherit the same bold/border/
# center+top style as the existing header cells instead of landing with no
# style at all.
$carSheet.Range("G1").Copy()
$carSheet.Range("H1:N1").PasteSpecial(-4122)

# Header row
$carSheet.Cells.Item(1, 2).Value = "name"
$carSheet.Cells.Item(1, 3).Value = "capacity"
$carSheet.Cells.Item(1, 4).Value = "owner"
$carSheet.Cells.Item(1, 5).Value = "register_date"
$carSheet.Cells.Item(1, 6).Value = "register_reason"
$carSheet.Cells.Item(1, 7).Value = "acquire_value"
$carSheet.Cells.Item(1, 8).Value = "property_category"
$carSheet.Cells.Item(1, 9).Value = "category"
$carSheet.Cells.Item(1, 10).Value = "date"
$carSheet.Cells.Item(1, 11).Value = "legislator_name"
$carSheet.Cells.Item(1, 12).Value = "legislator_id"
$carSheet.Cells.Item(1, 13).Value = "source_file"
$carSheet.Cells.Item(1, 14).Value = "index"

# Data row (row 2) - single, merged, correctly-shaped car record
$carSheet.Cells.Item(2, 1).Value = 33
$carSheet.Cells.Item(2, 2).Value = "三陽本田"
$carSheet.Cells.Item(2, 3).Value = 1997
$carSheet.Cells.Item(2, 4).Value = "廖靖汝"
$carSheet.Cells.Item(2, 5).Value = "97年05月06日"
$carSheet.Cells.Item(2, 6).Value = "買賣"
$carSheet.Cells.Item(2, 7).Value = 830000
$carSheet.Cells.Item(2, 8).Value = "land"
$carSheet.Cells.Item(2, 9).Value = "normal"
$carSheet.Cells.Item(2, 10).Value = "2013-12-30"
$carSheet.Cells.Item(2, 11).Value = "蘇震清"
$carSheet.Cells.Item(2, 12).Value = 1718
$carSheet.Cells.Item(2, 13).Value = "tmpb3b61"
$carSheet.Cells.Item(2, 14).Value = 33

# ---------------------------------------------------------------------------
# Sheets 4-7 ("存款", "具有相當價值之財產", "保險", "事業投資"): collapsing
# the car sheet down to a single data row shifts the running "index" column
# (column A) for every row that follows it by -3.
# ---------------------------------------------------------------------------
$depositSheet = $wb.Worksheets.Item(4)
$depositSheet.Cells.Item(2, 1).Value = 43
$depositSheet.Cells.Item(3, 1).Value = 44
$depositSheet.Cells.Item(4, 1).Value = 45
$depositSheet.Cells.Item(5, 1).Value = 46
$depositSheet.Cells.Item(6, 1).Value = 47
$depositSheet.Cells.Item(7, 1).Value = 48

$valuablesSheet = $wb.Worksheets.Item(5)
$valuablesSheet.Cells.Item(2, 1).Value = 75

$insuranceSheet = $wb.Worksheets.Item(6)
$insuranceSheet.Cells.Item(2, 1).Value = 80
$insuranceSheet.Cells.Item(3, 1).Value = 81
$insuranceSheet.Cells.Item(4, 1).Value = 82
$insuranceSheet.Cells.Item(5, 1).Value = 83
$insuranceSheet.Cells.Item(6, 1).Value = 84
$insuranceSheet.Cells.Item(7, 1).Value = 85

$investmentSheet = $wb.Worksheets.Item(7)
$investmentSheet.Cells.Item(2, 1).Value = 98
